# Rename the worksheet and refresh its contents to reflect the latest
# SUNAT processing run: the header row loses its bold/boxed styling and a
# fresh duplicate pull for the first ten companies is inserted, pushing
# the remaining rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet
$ws.Name = "Resultados"

# 2. Strip the bold + thin-border header formatting (A1:G1 goes back to
#    the workbook's default "Normal" style)
$ws.Range("A1:G1").ClearFormats()

# 3. Insert 10 fresh rows right after row 11 and repopulate them with a
#    copy of the first data block (rows 2-11), shifting the remaining
#    companies (old rows 12-18) down to rows 22-28
$ws.Rows("12:21").Insert()
$ws.Range("A2:G11").Copy()
$ws.Range("A12").PasteSpecial()
